$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 400; existing rows 400-460 shift down to 401-461.
$ws.Rows(400).Insert()

# Populate the newly inserted row 400 with the new weekly price record.
$ws.Range("A400").Value = 4
$ws.Range("B400").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C400").Value = "Los Lagos"
$ws.Range("D400").Value = 45034
$ws.Range("E400").Value = 10
$ws.Range("F400").Value = 100112045
$ws.Range("G400").Value = "Zapallo"
$ws.Range("H400").Value = "Paine"
$ws.Range("I400").Value = "1a (cosecha)"
$ws.Range("J400").Value = 1200
$ws.Range("K400").Value = 580
$ws.Range("L400").Value = 600
$ws.Range("M400").Value = 590
$ws.Range("N400").Value = "$/kilo (volumen en unidades)"
$ws.Range("O400").Value = "Región de O'Higgins"
$ws.Range("P400").Value = 590
$ws.Range("Q400").Value = 1
$ws.Range("R400").Value = "Hortaliza"
